$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Row 2 -----
$ws.Range("D2").Value = -0.012
$ws.Range("G2").Value = -0.008859504132231404
$ws.Range("H2").Value = -0.008859504132231404
$ws.Range("I2").Value = -0.02987549812133724
$ws.Range("J2").Value = -0.02987549812133724
$ws.Range("K2").Value = -2.307
$ws.Range("L2").Value = -0.02542148760330578
$ws.Range("U2").Value = 6.49
$ws.Range("V2").Value = 0.2438932732055618
$ws.Range("X2").Value = 0.06016853439930246
$ws.Range("Z2").Value = -5.26973104491134
$ws.Range("AA2").Value = 0.4193476633270933
$ws.Range("AB2").Value = 0.05244901121731119
$ws.Range("AC2").Value = 0.3668986521097821
$ws.Range("AD2").Value = 9.199999999999999
$ws.Range("AE2").Value = 0.5610072725567721
$ws.Range("AF2").Value = 9.761007272556771
$ws.Range("AG2").Value = 3.271007272556771
$ws.Range("AH2").Value = 0.2683733007285119
$ws.Range("AI2").Value = 0.8569046651451977
$ws.Range("AJ2").Value = 0.1094677713746591
$ws.Range("AK2").Value = 0.6674153068233141
$ws.Range("AL2").Value = 1.06
$ws.Range("AM2").Value = 0.9760000000000001
$ws.Range("AN2").Value = -6.764705882352941
$ws.Range("AO2").Value = -2.70377358490566
$ws.Range("AP2").Value = -2.405152406291744
$ws.Range("AQ2").Value = -2.936475409836065
$ws.Range("W2").ClearContents()
$ws.Range("Y2").ClearContents()

# ----- Row 3 -----
$ws.Range("B3").Value = "DCI Database for Commerce and Industry AG (DB:DCIK)"
$ws.Range("G3").Value = 0.01885714285714286
$ws.Range("H3").Value = 0.01885714285714286
$ws.Range("I3").Value = 0.04676812170085377
$ws.Range("J3").Value = 0.04676812170085377
$ws.Range("K3").Value = 0.113
$ws.Range("L3").Value = 0.06457142857142857
$ws.Range("O3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("X3").Value = 0.05454743392931438
$ws.Range("Z3").Value = 14.48928158123833
$ws.Range("AA3").Value = 0.6776364843492929
$ws.Range("AB3").Value = 0.05290927020913389
$ws.Range("AC3").Value = 0.624727214140159
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0.1207789351175296
$ws.Range("AF3").Value = 0.1207789351175296
$ws.Range("AG3").Value = 0.1207789351175296
$ws.Range("AH3").Value = 0.03985079007844098
$ws.Range("AI3").Value = 1
$ws.Range("AJ3").Value = 0.03985079007844098
$ws.Range("AK3").Value = 1
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
$ws.Range("AN3").Value = 0
$ws.Range("AP3").Value = 1.050251609717648
$ws.Range("D3").ClearContents()
$ws.Range("W3").ClearContents()
$ws.Range("Y3").ClearContents()
$ws.Range("AO3").ClearContents()
$ws.Range("AQ3").ClearContents()

# ----- Row 4 -----
$ws.Range("B4").Value = "asknet Solutions AG (XTRA:ASKN)"
$ws.Range("D4").Value = -0.012
$ws.Range("G4").Value = -0.009404494382022472
$ws.Range("H4").Value = -0.009404494382022472
$ws.Range("I4").Value = -0.03138253558975111
$ws.Range("J4").Value = -0.03138253558975111
$ws.Range("K4").Value = -2.42
$ws.Range("L4").Value = -0.02719101123595506
$ws.Range("O4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("U4").Value = 6.49
$ws.Range("V4").Value = 0.2738396624472574
$ws.Range("W4").Value = 2.515592515592516
$ws.Range("X4").Value = 0.06578963486929054
$ws.Range("Y4").Value = 2.449802880723225
$ws.Range("Z4").Value = -5.132116933135647
$ws.Range("AA4").Value = 0.1610588423048938
$ws.Range("AB4").Value = 0.05198875222548849
$ws.Range("AC4").Value = 0.1090700900794053
$ws.Range("AD4").Value = 9.199999999999999
$ws.Range("AE4").Value = 0.4402283374392424
$ws.Range("AF4").Value = 9.640228337439241
$ws.Range("AG4").Value = 3.150228337439241
$ws.Range("AH4").Value = 0.2891470400223323
$ws.Range("AI4").Value = 0.8553711645233303
$ws.Range("AJ4").Value = 0.117325942180039
$ws.Range("AK4").Value = 0.6590121046658647
$ws.Range("AL4").Value = 1.06
$ws.Range("AM4").Value = 0.9760000000000001
$ws.Range("AN4").Value = -6.237288135593221
$ws.Range("AO4").Value = -2.726415094339623
$ws.Range("AP4").Value = -2.135748025382536
$ws.Range("AQ4").Value = -2.961065573770492
